# Apply the "TODO CMS" sheet population + selection/active-tab changes
# described by the commit: "Products are shown now on the CMS-homepage.
# Deleting and adding products is blind coded in backend, not tested yet."

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # TODO
$ws2 = $wb.Worksheets.Item(2)   # TODO CMS
$ws3 = $wb.Worksheets.Item(3)   # Testfälle

# ---------------------------------------------------------------------
# 1. Populate the "TODO CMS" sheet (sheet2) with a TODO table identical
#    in shape to the "TODO" sheet's table.
# ---------------------------------------------------------------------

# Header row (copy the bold header style from the TODO sheet)
$ws2.Range("A1").Value = "TODO"
$ws2.Range("B1").Value = "Status"
$ws2.Range("C1").Value = "Name"
$ws2.Range("D1").Value = "Hinweis"
$ws1.Range("A1:D1").Copy()
$ws2.Range("A1:D1").PasteSpecial(-4122)

# Data rows
$rows = @(
    @{ A = "Löschen von Bestellungen";            B = "offen"; C = "Jonas" },
    @{ A = "Logout";                               B = "offen"; C = "Jonas" },
    @{ A = "Löschen von Usern";                    B = "offen"; C = "Jonas" },
    @{ A = "Produkte löschen";                     B = "offen"; C = "Jonas" },
    @{ A = "Produkte hinzufügen";                  B = "offen"; C = "Jonas" },
    @{ A = "Verifikation von neuen Admin-Usern";   B = "offen"; C = "Jonas" },
    @{ A = "Passwort ändern";                      B = "offen"; C = "Jonas" },
    @{ A = "Anzeigen aller Produkte";               B = "done";  C = "Jonas" }
)

$r = 2
foreach ($row in $rows) {
    $ws2.Range("A$r").Value = $row.A
    $ws2.Range("B$r").Value = $row.B
    $ws2.Range("C$r").Value = $row.C

    if ($row.B -eq "offen") {
        $ws1.Range("B4").Copy()
    } else {
        $ws1.Range("B2").Copy()
    }
    $ws2.Range("B$r").PasteSpecial(-4122)

    $r = $r + 1
}

# Column A width on the CMS sheet
$ws2.Columns.Item(1).ColumnWidth = 36

# ---------------------------------------------------------------------
# 2. Selections on each sheet.
# ---------------------------------------------------------------------
$ws1.Range("B6").Select()
$ws3.Range("C5").Select()

# TODO CMS becomes the active/selected tab
$ws2.Activate()
$ws2.Range("C16").Select()

# ---------------------------------------------------------------------
# 3. Add an "OK" result to the 3rd Testfall row (Testfälle sheet).
# ---------------------------------------------------------------------
$ws3.Range("B3").Value = "OK"
$ws3.Range("B2").Copy()
$ws3.Range("B3").PasteSpecial(-4122)
